$d = $word.ActiveDocument

$d.Content.Find.Execute("59826 (12.3%)", $true, $false, $false, $false, $false, $true, 1, $false, "59826 (26.3%)", 2)
$d.Content.Find.Execute("52790 (10.8%)", $true, $false, $false, $false, $false, $true, 1, $false, "52790 (20.4%)", 2)
$d.Content.Find.Execute("43037 (8.8%)", $true, $false, $false, $false, $false, $true, 1, $false, "43037 (18.9%)", 2)
$d.Content.Find.Execute("61962 (12.7%)", $true, $false, $false, $false, $false, $true, 1, $false, "61962 (23.9%)", 2)
$d.Content.Find.Execute("124525 (25.6%)", $true, $false, $false, $false, $false, $true, 1, $false, "124525 (54.8%)", 2)
$d.Content.Find.Execute("144487 (29.7%)", $true, $false, $false, $false, $false, $true, 1, $false, "144487 (55.7%)", 2)
